$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for account 005135532 (FELIPE) entirely
$felipeCell = $ws.Columns(1).Find("005135532")
$felipeCell.EntireRow.Delete()

# Delete the row for account 004589191 (CONTEL) entirely
$contelCell = $ws.Columns(1).Find("004589191")
$contelCell.EntireRow.Delete()

# Update MARCELO's (004748761) balance
$marceloCell = $ws.Columns(1).Find("004748761")
$ws.Cells.Item($marceloCell.Row, 3).Value = 29026.7

# Update ROSANGELA's (005002457) balance
$rosangelaCell = $ws.Columns(1).Find("005002457")
$ws.Cells.Item($rosangelaCell.Row, 3).Value = 3000

# Re-sort the data rows (account, name, balance) descending by balance,
# anchored between the header row and the last account row (004472404 / DILSON),
# leaving the trailing blank row and the filter-summary row untouched.
$lastDataCell = $ws.Columns(1).Find("004472404")
$lastRow = $lastDataCell.Row

$dataRange = $ws.Range("A2:C" + $lastRow)
$sortKey = $ws.Range("C2:C" + $lastRow)
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, $null, 2)
